$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sale report gained a new line item ("معجون سيجنال 50 مل") as row #13.
# It is inserted right before the "total" row, which (together with the
# footer row below it) shifts down by one row. Concretely:
#   old row 19 (total)   -> new row 20
#   old row 20 (footer)  -> new row 21
#   new row 19           -> the new item row
# ---------------------------------------------------------------------------

# Insert a blank row at position 19; this pushes the total/footer rows down.
$ws.Rows.Item(19).Insert()

# Clone the formatting (styles) of the last item row (18) onto the new row 19
# so every cell in the new row carries the same look as the other item rows.
$ws.Range("A18:Q18").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merged cells that make up one item row.
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# Match the exact row heights seen in the final report.
$ws.Rows.Item(19).RowHeight = 25.5
$ws.Rows.Item(20).RowHeight = 24.75

# ---------------------------------------------------------------------------
# Fill in the new item's data (row #13 in the listing).
# ---------------------------------------------------------------------------
$ws.Range("A19").Value() = 13
$ws.Range("C19").Value() = "معجون سيجنال 50 مل"
$ws.Range("H19").Value() = "3:0"

# L19 / P19 sit on cells whose number format is numeric, but the source
# workbook always stores these columns as literal text - temporarily switch
# to a text format so the assignment isn't re-interpreted as a number.
$fmtL = $ws.Range("L19").NumberFormat
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value() = "0"
$ws.Range("L19").NumberFormat = $fmtL

$ws.Range("N19").Value() = "35.00"

$fmtP = $ws.Range("P19").NumberFormat
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value() = "35.0000"
$ws.Range("P19").NumberFormat = $fmtP

$ws.Range("Q19").Value() = "1:0"

# ---------------------------------------------------------------------------
# Update the running total (now on row 20) to include the new item's price.
# ---------------------------------------------------------------------------
$oldTotal = $ws.Range("P20").Value()
$ws.Range("P20").Value() = $oldTotal + 35

# ---------------------------------------------------------------------------
# The footer timestamp (now on row 21) was refreshed to the new save time.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value() = "Tuesday, 29 July, 2025 9:45 AM"
